# Add 2022-Q4 data:
#  - insert a new summary row into "总计" (total) sheet
#  - insert a new "2022-Q4" sheet (positioned before "2022-Q3") with the
#    quarterly fund-holding detail

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for the 2022-Q4 summary line,
#    pushing the existing 2022-Q3 / 2022-Q2 / 2022-Q1 / 2021-Q4 rows down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows("2:2").Insert()

# The Insert() operation drags the row-above formatting onto the new row;
# reset it back to the plain (unstyled) look the sibling data rows use,
# then re-apply the "A" column style from an untouched sibling row (A4,
# still holding its original formatting at this point).
$total.Range("B2:D2").ClearFormats()
$total.Range("A4").Copy()
$total.Range("A2").PasteSpecial($xlPasteFormats)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.53

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, duplicated from "2022-Q3" so it inherits the
#    exact same layout/formatting, then positioned right before it.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Clean, never-overwritten reference cells on row 2 used purely as
# formatting donors (PasteSpecial formats only) so every edited/added
# cell ends up with the same "no explicit style" look the source data
# uses, instead of inheriting the "@" text-number-format we set below
# in order to stop Excel from renumbering things like "015870" -> 15870.
$cleanStyleRef = "A1"          # default/no style (used as a safety net, unused directly)
$aStyleRef = "A2"              # s="2" style (index/"A" column look)
$plainStyleRef = "B2"          # no style, never itself modified

function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($plainStyleRef).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

function Set-IndexCell($ws, $addr, $num) {
    $ws.Range($aStyleRef).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
    $ws.Range($addr).Value = $num
}

# Row 2: update the existing fund's figures for the new quarter
# (B2/C2 - code & name - stay the same: 540002 / 汇丰晋信龙腾混合).
Set-TextCell $q4 "D2" "9.82"
Set-TextCell $q4 "E2" "90.84"
Set-TextCell $q4 "F2" "5.32"
Set-TextCell $q4 "G2" "0.5224"
$q4.Range("H2").Value = 1

# Row 3: 015870 富安达先进制造混合A
Set-IndexCell $q4 "A3" 1
Set-TextCell $q4 "B3" "015870"
Set-TextCell $q4 "C3" "富安达先进制造混合A"
Set-TextCell $q4 "D3" "0.10"
Set-TextCell $q4 "E3" "83.11"
Set-TextCell $q4 "F3" "3.46"
Set-TextCell $q4 "G3" "0.0035"
$q4.Range("H3").Value = 2

# Row 4: 001797 国新国证新利灵活配置混合
Set-IndexCell $q4 "A4" 2
Set-TextCell $q4 "B4" "001797"
Set-TextCell $q4 "C4" "国新国证新利灵活配置混合"
Set-TextCell $q4 "D4" "0.02"
Set-TextCell $q4 "E4" "81.37"
Set-TextCell $q4 "F4" "2.93"
Set-TextCell $q4 "G4" "0.0006"
$q4.Range("H4").Value = 10

# Row 5: 015886 富安达先进制造混合C (holding-value column G is a bare 0,
# stored as a real number rather than text, per the source data).
Set-IndexCell $q4 "A5" 3
Set-TextCell $q4 "B5" "015886"
Set-TextCell $q4 "C5" "富安达先进制造混合C"
Set-TextCell $q4 "D5" "0.00"
Set-TextCell $q4 "E5" "83.11"
Set-TextCell $q4 "F5" "3.46"
$q4.Range("G5").Value = 0
$q4.Range("H5").Value = 2
